# Update the "as_of_utc" timestamp column (AA, rows 2-26) on the
# "Главные" and "Линейные" sheets from 2025-12-17 04:00:08 to
# 2025-12-17 07:03:58 (republish timestamp refresh).

$wb = $excel.ActiveWorkbook

$oldTs = "2025-12-17 04:00:08"
$newTs = "2025-12-17 07:03:58"

$sheetNames = @("Главные", "Линейные")

foreach ($sheetName in $sheetNames) {
    $ws = $wb.Worksheets.Item($sheetName)
    for ($row = 2; $row -le 26; $row++) {
        $cell = $ws.Cells.Item($row, 27)  # column AA = 27
        if ($cell.Value2 -eq $oldTs) {
            $cell.Value = $newTs
        }
    }
}
